$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: A5 and B5 get plain/default style, C5 keeps its existing style (index 5)
$ws.Range("A5").Value = "Leo"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "LCM San Pedro"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "Leo.xlsx"

# Row 6: all cells get the default/plain style
$ws.Range("A6").Value = "Instituto Costarricense de Electricidad"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "Sabana Norte edificio ICE"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "Instituto Costarricense de Electricidad.xlsx"
$ws.Range("C6").Style = "Normal"
